# Update "想去人数" (interest count) figures on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 8206
$ws1.Range("F9").Value = 139
$ws1.Range("F10").Value = 210
$ws1.Range("F14").Value = 5013
$ws1.Range("F16").Value = 5463
$ws1.Range("F19").Value = 344
$ws1.Range("F22").Value = 255
$ws1.Range("F28").Value = 1676
$ws1.Range("F29").Value = 839
$ws1.Range("F40").Value = 4811
$ws1.Range("F44").Value = 149
$ws1.Range("F45").Value = 936
$ws1.Range("F46").Value = 892
$ws1.Range("F50").Value = 959

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F10").Value = 8206
$ws4.Range("F11").Value = 139
$ws4.Range("F12").Value = 210
$ws4.Range("F17").Value = 5013
$ws4.Range("F19").Value = 5463
$ws4.Range("F22").Value = 344
$ws4.Range("F25").Value = 255
$ws4.Range("F29").Value = 1676
$ws4.Range("F30").Value = 839
$ws4.Range("F41").Value = 4811
$ws4.Range("F45").Value = 149
$ws4.Range("F46").Value = 938
$ws4.Range("F50").Value = 962
